$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Append a new data row (row 46) to each of the four log sheets, continuing
# the existing time-series with the next sample captured on 2025-06-24.
# ---------------------------------------------------------------------------

# --- FE_LFT_#1 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("FE_LFT_#1")
$ws.Cells.Item(46, 1).Value = 45832.49384259259
$ws.Cells.Item(46, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 2).Value = "0x01,0x7c"
$ws.Cells.Item(46, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(46, 4).Value = "0x01,0x5C"
$ws.Cells.Item(46, 5).Value = "0xf"
$ws.Cells.Item(46, 6).Value = 380
$ws.Cells.Item(46, 7).Value = 759863127514710900000000.0
$ws.Cells.Item(46, 8).Value = 348
$ws.Cells.Item(46, 9).Value = 15

# --- FE_LFT_#2 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("FE_LFT_#2")
$ws.Cells.Item(46, 1).Value = 45832.49384259259
$ws.Cells.Item(46, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 2).Value = "0x01,0x90"
$ws.Cells.Item(46, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(46, 4).Value = "0x01,0x70"
$ws.Cells.Item(46, 5).Value = "0xe"
$ws.Cells.Item(46, 6).Value = 400
$ws.Cells.Item(46, 7).Value = 568432987514711000000000.0
$ws.Cells.Item(46, 8).Value = 368
$ws.Cells.Item(46, 9).Value = 14

# --- FE_PLT_#1 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("FE_PLT_#1")
$ws.Cells.Item(46, 1).Value = 45832.49384259259
$ws.Cells.Item(46, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 2).Value = "0x00,0x6e"
$ws.Cells.Item(46, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(46, 4).Value = "0x00,0x69"
$ws.Cells.Item(46, 5).Value = "0x3"
$ws.Cells.Item(46, 6).Value = 110
$ws.Cells.Item(46, 7).Value = 568631262647114000000000.0
$ws.Cells.Item(46, 8).Value = 105
$ws.Cells.Item(46, 9).Value = 3

# --- FE_PLT_#2 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("FE_PLT_#2")
$ws.Cells.Item(46, 1).Value = 45832.49384259259
$ws.Cells.Item(46, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 2).Value = "0x00,0x6e"
$ws.Cells.Item(46, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(46, 4).Value = "0x00,0x69"
$ws.Cells.Item(46, 5).Value = "0x3"
$ws.Cells.Item(46, 6).Value = 110
$ws.Cells.Item(46, 7).Value = 985046333984776000000000.0
$ws.Cells.Item(46, 8).Value = 105
$ws.Cells.Item(46, 9).Value = 3
